$d = $word.ActiveDocument

# Locate the paragraph that contains the "missing expression" AQL error
# message (the second paragraph of the document in this fixture).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Expression*self.name*is invalid: missing expression*") {
        $target = $p
        break
    }
}

# Collapsed range positioned right before the paragraph mark, i.e. right
# after the last existing run of the paragraph.
$rng = $d.Range($target.Range.End - 1, $target.Range.End - 1)

function Add-PlainRun($range, $text) {
    $range.InsertAfter($text)
    $range.Collapse(0)
}

function Add-ErrorRun($range, $text) {
    $start = $range.End
    $range.InsertAfter($text)
    $newRng = $d.Range($start, $range.End)
    $newRng.Font.Color = 255
    $newRng.Font.Size = 16
    $newRng.Font.HighlightColorIndex = 16
    $range.Collapse(0)
}

Add-PlainRun $rng "    "
Add-ErrorRun $rng "<---"
Add-ErrorRun $rng "Couldn't find the 'self' variable"

Add-PlainRun $rng "    "
Add-ErrorRun $rng "<---"
Add-ErrorRun $rng "missing expression"

Add-PlainRun $rng "    "
Add-ErrorRun $rng "<---"
Add-ErrorRun $rng "The predicate never evaluates to a boolean type ([])."

$d.Saved = $false
